$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 265, shifting existing rows 265-294 down to 267-296
$ws.Rows.Item(265).Insert()
$ws.Rows.Item(265).Insert()

# --- New row 265 ---
$ws.Range("A265").Value = 4
$ws.Range("B265").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C265").Value = "Los Lagos"
$ws.Range("D265").Value = 44491
$ws.Range("E265").Value = 10
$ws.Range("F265").Value = "Fruta"
$ws.Range("G265").Value = 100102
$ws.Range("H265").Value = "Cítricos"
$ws.Range("I265").Value = 100102003
$ws.Range("J265").Value = "Limón"
$ws.Range("K265").Value = "Sin especificar"
$ws.Range("L265").Value = "1a amarillo"
$ws.Range("M265").Value = 1200
$ws.Range("N265").Value = 9500
$ws.Range("O265").Value = 10000
$ws.Range("P265").Value = 9750
$ws.Range("Q265").Value = "`$/malla 18 kilos"
$ws.Range("R265").Value = "Provincia de Melipilla"
$ws.Range("S265").Value = 542
$ws.Range("T265").Value = 18

# --- New row 266 ---
$ws.Range("A266").Value = 4
$ws.Range("B266").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C266").Value = "Los Lagos"
$ws.Range("D266").Value = 44491
$ws.Range("E266").Value = 10
$ws.Range("F266").Value = "Fruta"
$ws.Range("G266").Value = 100102
$ws.Range("H266").Value = "Cítricos"
$ws.Range("I266").Value = 100102003
$ws.Range("J266").Value = "Limón"
$ws.Range("K266").Value = "Sin especificar"
$ws.Range("L266").Value = "2a amarillo"
$ws.Range("M266").Value = 500
$ws.Range("N266").Value = 8000
$ws.Range("O266").Value = 8000
$ws.Range("P266").Value = 8000
$ws.Range("Q266").Value = "`$/malla 18 kilos"
$ws.Range("R266").Value = "Provincia de Melipilla"
$ws.Range("S266").Value = 444
$ws.Range("T266").Value = 18

# Ensure D column (date) keeps the same date number format used by the rest
# of the column, matching the style applied to the surrounding rows.
$ws.Range("D265").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D266").NumberFormat = "YYYY-MM-DD HH:MM:SS"
